# New crime data collected — refresh the CompStat weekly report:
#  - bump the report "Volume/Number" label and the week-covering dates
#  - refresh every stat cell in the precinct/subdivision table (rows 15-29)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text (shared strings with uniform run formatting, so writing
# the full, concatenated label back through .Value reproduces the same
# visible text the diff produces).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/25/2023  Through  12/31/2023"

# ---------------------------------------------------------------------
# Helper: some "C" cells flip between a real number and the literal
# text "0" (used by this report to mean "no prior value"/blank stat).
# Style must track the sibling "D" cell of the same row, so copy
# formats across via PasteSpecial before/after setting the value.
# ---------------------------------------------------------------------
function Set-TextZero($ref, $donorRef) {
    $ws.Range($ref).Value = "'0"
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

function Set-NumberFromDonor($ref, $donorRef, $value) {
    $ws.Range($donorRef).Copy() | Out-Null
    $ws.Range($ref).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Range($ref).Value2 = $value
}

# Row 15: C15 numeric 1 -> text "0"
Set-TextZero "C15" "D15"
$ws.Range("L15").Value2 = 11.764705882352

# Row 16
$ws.Range("C16").Value2 = 1
$ws.Range("D16").Value2 = 1
$ws.Range("F16").Value2 = 10
$ws.Range("G16").Value2 = 9
$ws.Range("H16").Value2 = 11.111111111111
$ws.Range("I16").Value2 = 113
$ws.Range("J16").Value2 = 110
$ws.Range("K16").Value2 = 2.727272727272
$ws.Range("L16").Value2 = 88.333333333333
$ws.Range("M16").Value2 = -50.869565217391
$ws.Range("N16").Value2 = -84.129213483146

# Row 17
$ws.Range("C17").Value2 = 4
$ws.Range("D17").Value2 = 6
$ws.Range("E17").Value2 = -33.333333333333
$ws.Range("F17").Value2 = 18
$ws.Range("G17").Value2 = 15
$ws.Range("H17").Value2 = 20
$ws.Range("I17").Value2 = 239
$ws.Range("J17").Value2 = 216
$ws.Range("K17").Value2 = 10.648148148148
$ws.Range("L17").Value2 = 10.138248847926
$ws.Range("M17").Value2 = 49.375
$ws.Range("N17").Value2 = -13.090909090909

# Row 18
$ws.Range("C18").Value2 = 2
$ws.Range("G18").Value2 = 8
$ws.Range("H18").Value2 = -50
$ws.Range("I18").Value2 = 55
$ws.Range("K18").Value2 = -42.105263157894
$ws.Range("L18").Value2 = -27.631578947368
$ws.Range("M18").Value2 = -78.260869565217
$ws.Range("N18").Value2 = -91.40625

# Row 19
$ws.Range("C19").Value2 = 5
$ws.Range("D19").Value2 = 7
$ws.Range("E19").Value2 = -28.571428571428
$ws.Range("F19").Value2 = 18
$ws.Range("G19").Value2 = 28
$ws.Range("H19").Value2 = -35.714285714285
$ws.Range("I19").Value2 = 219
$ws.Range("J19").Value2 = 247
$ws.Range("K19").Value2 = -11.336032388664
$ws.Range("L19").Value2 = 15.873015873015
$ws.Range("M19").Value2 = -12.048192771084
$ws.Range("N19").Value2 = -47.101449275362

# Row 20
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = 2
$ws.Range("E20").Value2 = -50
$ws.Range("F20").Value2 = 17
$ws.Range("G20").Value2 = 10
$ws.Range("H20").Value2 = 70
$ws.Range("I20").Value2 = 166
$ws.Range("J20").Value2 = 153
$ws.Range("K20").Value2 = 8.496732026143
$ws.Range("L20").Value2 = 72.916666666666
$ws.Range("M20").Value2 = 48.214285714285
$ws.Range("N20").Value2 = -92.378328741965

# Row 21
$ws.Range("C21").Value2 = 13
$ws.Range("D21").Value2 = 16
$ws.Range("E21").Value2 = -18.75
$ws.Range("F21").Value2 = 69
$ws.Range("G21").Value2 = 70
$ws.Range("H21").Value2 = -1.428571428571
$ws.Range("I21").Value2 = 815
$ws.Range("J21").Value2 = 839
$ws.Range("K21").Value2 = -2.860548271752
$ws.Range("L21").Value2 = 22.372372372372
$ws.Range("M21").Value2 = -20.719844357976
$ws.Range("N21").Value2 = -80.837056195626

# Row 22
$ws.Range("M22").Value2 = -23.076923076923

# Row 23: C23 text "0" -> numeric 4
Set-NumberFromDonor "C23" "D23" 4
$ws.Range("D23").Value2 = 3
$ws.Range("E23").Value2 = 33.333333333333
$ws.Range("G23").Value2 = 13
$ws.Range("H23").Value2 = -15.384615384615
$ws.Range("I23").Value2 = 127
$ws.Range("J23").Value2 = 107
$ws.Range("K23").Value2 = 18.691588785046
$ws.Range("L23").Value2 = 38.043478260869
$ws.Range("M23").Value2 = 111.666666666667

# Row 24
$ws.Range("C24").Value2 = 9
$ws.Range("D24").Value2 = 21
$ws.Range("E24").Value2 = -57.142857142857
$ws.Range("F24").Value2 = 39
$ws.Range("G24").Value2 = 49
$ws.Range("H24").Value2 = -20.408163265306
$ws.Range("I24").Value2 = 542
$ws.Range("J24").Value2 = 567
$ws.Range("K24").Value2 = -4.409171075837
$ws.Range("L24").Value2 = 22.072072072072
$ws.Range("M24").Value2 = 0.743494423791

# Row 25
$ws.Range("C25").Value2 = 6
$ws.Range("D25").Value2 = 5
$ws.Range("E25").Value2 = 20
$ws.Range("G25").Value2 = 29
$ws.Range("H25").Value2 = 31.034482758620
$ws.Range("I25").Value2 = 378
$ws.Range("J25").Value2 = 383
$ws.Range("K25").Value2 = -1.305483028720
$ws.Range("L25").Value2 = 18.867924528301
$ws.Range("M25").Value2 = -16.556291390728

# Row 26: C26 numeric 2 -> text "0"
Set-TextZero "C26" "D26"
$ws.Range("L26").Value2 = -7.692307692307

# Row 27
$ws.Range("G27").Value2 = 1

# Row 28: C28 numeric 1 -> text "0"
Set-TextZero "C28" "D28"

# Row 29: C29 numeric 1 -> text "0"
Set-TextZero "C29" "D29"
